# US9430_02_TransactionHistory: ITG2 -> ITG4
# Update the quote number used by the transaction-history test data and
# leave the selection where the author left it before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A1 already holds the header "strQuoteNumber"; A2 holds the quote number
# value itself - bump it to the new ITG4 quote number.
$ws.Range("A2").Value = "NI00164529"

# Move the active selection from B10 to D10, matching the author's cursor
# position when the sheet was last saved.
$ws.Range("D10").Select()
